$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark (it will be re-added on the
#     new paragraph we are about to insert). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: find the "Drop any unwanted columns" paragraph and insert a
#     brand-new paragraph right before it: a tab followed by
#     "Format all dates correctly for openrefine", with a _GoBack bookmark
#     positioned right after the text (inside the same paragraph). ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Drop any unwanted columns*") {
        $target = $p
        break
    }
}

$newPara = $target.Range.InsertParagraphBefore()

# Re-fetch the freshly inserted (still empty) paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Drop any unwanted columns*") {
        $target = $p
        break
    }
}
$inserted = $target.Previous()
$r = $inserted.Range

# Build the run (tab + text) via a raw WordML fragment so the tab becomes a
# real <w:tab/> element rather than a literal tab character, and append a
# trailing placeholder character "X" that we will use as a safe anchor for
# the bookmark before removing it again.
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Format all dates correctly for openrefineX</w:t></w:r></w:p>'
$r.InsertXML($xmlFrag) | Out-Null

# Re-fetch the paragraph again now that it holds real content.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Drop any unwanted columns*") {
        $target = $p
        break
    }
}
$inserted = $target.Previous()
$r2 = $inserted.Range

# Place the bookmark immediately before the trailing "X" placeholder (i.e.
# right after the real text) -- this is a safe, unambiguous mid-paragraph
# position, unlike the paragraph-end boundary.
$bmPos = $r2.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove the "X" placeholder character now that the bookmark is anchored.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Drop any unwanted columns*") {
        $target = $p
        break
    }
}
$inserted = $target.Previous()
$r3 = $inserted.Range
$xRange = $d.Range($r3.End - 2, $r3.End - 1)
$xRange.Delete()
